$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C21").NumberFormat = "@"

$ws.Range("C2").Value = "111010100111001111000000111101"
$ws.Range("D2").Value = 0.8387423222178733
$ws.Range("E2").Value = 0.03472310536680396
$ws.Range("F2").Value = 0.3853178340835093

$ws.Range("C3").Value = "111010100111001111000000111101"
$ws.Range("D3").Value = 0.8387423222178733
$ws.Range("E3").Value = 0.02955388167217791
$ws.Range("F3").Value = 0.682852610050532

$ws.Range("C4").Value = "111010100111001111001000100110"
$ws.Range("D4").Value = 0.8387431563855785
$ws.Range("E4").Value = 0.824493532946885
$ws.Range("F4").Value = 0.8373164731692107

$ws.Range("C5").Value = "111010100111001111001000100110"
$ws.Range("D5").Value = 0.8387431563855785
$ws.Range("E5").Value = 0.838731841422177
$ws.Range("F5").Value = 0.838740392892275

$ws.Range("C6").Value = "111010100111001111001000111101"
$ws.Range("D6").Value = 0.8387431956204705
$ws.Range("E6").Value = 0.8387422829830017
$ws.Range("F6").Value = 0.8387429895520395

$ws.Range("C7").Value = "111010100111001111001000111101"
$ws.Range("D7").Value = 0.8387431956204705
$ws.Range("E7").Value = 0.8387422829830017
$ws.Range("F7").Value = 0.8387430886627669

$ws.Range("C8").Value = "111010100111001111001000111101"
$ws.Range("D8").Value = 0.8387431956204705
$ws.Range("E8").Value = 0.8387431563855785
$ws.Range("F8").Value = 0.8387431877734922

$ws.Range("C9").Value = "111010100111001111001000111101"
$ws.Range("D9").Value = 0.8387431956204705
$ws.Range("E9").Value = 0.8244943988989016
$ws.Range("F9").Value = 0.8373183159483135

$ws.Range("C10").Value = "111010110111001111001000111101"
$ws.Range("D10").Value = 0.8459133703485278
$ws.Range("E10").Value = 0.8387431956204705
$ws.Range("F10").Value = 0.8394602130932762

$ws.Range("C11").Value = "111010110111001111001000111101"
$ws.Range("D11").Value = 0.8459133703485278
$ws.Range("E11").Value = 0.8386872987710199
$ws.Range("F11").Value = 0.842322693299554

$ws.Range("C12").Value = "111010110111001111001000111101"
$ws.Range("D12").Value = 0.8459133703485278
$ws.Range("E12").Value = 0.4485456253674628
$ws.Range("F12").Value = 0.8054595783776156

$ws.Range("C13").Value = "111010110111001111001000111101"
$ws.Range("D13").Value = 0.8459133703485278
$ws.Range("E13").Value = 0.8459133703485278
$ws.Range("F13").Value = 0.8459133703485279

$ws.Range("C14").Value = "111010110111001111001000111101"
$ws.Range("D14").Value = 0.8459133703485278
$ws.Range("E14").Value = 0.8459133703485278
$ws.Range("F14").Value = 0.8459133703485279

$ws.Range("C15").Value = "111010110111001111001000111101"
$ws.Range("D15").Value = 0.8459133703485278
$ws.Range("E15").Value = 0.8459133703485278
$ws.Range("F15").Value = 0.8459133703485279

$ws.Range("C16").Value = "111010110111001111001000111101"
$ws.Range("D16").Value = 0.8459133703485278
$ws.Range("E16").Value = 0.8459133703485278
$ws.Range("F16").Value = 0.8459133703485279

$ws.Range("C17").Value = "111010110111001111001100111101"
$ws.Range("D17").Value = 0.845913808912642
$ws.Range("E17").Value = 0.8459133703485278
$ws.Range("F17").Value = 0.8459134142049394

$ws.Range("C18").Value = "111010110111001111001100111101"
$ws.Range("D18").Value = 0.845913808912642
$ws.Range("E18").Value = 0.8459133703485278
$ws.Range("F18").Value = 0.8459135019177622

$ws.Range("C19").Value = "111010110111001111001100111101"
$ws.Range("D19").Value = 0.845913808912642
$ws.Range("E19").Value = 0.8459133703485278
$ws.Range("F19").Value = 0.8459137650562306

$ws.Range("C20").Value = "111011110111001111001100111101"
$ws.Range("D20").Value = 0.8748996910572698
$ws.Range("E20").Value = 0.845913808912642
$ws.Range("F20").Value = 0.8488123971271048

$ws.Range("C21").Value = "111011110111001111001100111101"
$ws.Range("D21").Value = 0.8748996910572698
$ws.Range("E21").Value = 0.845913808912642
$ws.Range("F21").Value = 0.8633053381994188
